$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics (per commit: "update scripts wuth new tpm")

# Row 2
$ws.Range("G2").Value = 35.73885133333334
$ws.Range("H2").Value = 107.216554
$ws.Range("I2").Value = 0.01949729408921566
$ws.Range("J2").Value = 0.01949729408921566
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8180823333333334
$ws.Range("N2").Value = 2.454247
$ws.Range("O2").Value = 0.5115352725808422
$ws.Range("P2").Value = 0.5115352725808422
$ws.Range("Q2").Value = 29.23732288942645
$ws.Range("R2").Value = 263.135906004838
$ws.Range("S2").Value = 0.009973553646515776
$ws.Range("T2").Value = 0.009973553646515775

# Row 3
$ws.Range("G3").Value = 35.73885133333334
$ws.Range("H3").Value = 107.216554
$ws.Range("I3").Value = 0.01949729408921566
$ws.Range("J3").Value = 0.01949729408921566
$ws.Range("M3").Value = 0.7811863333333333
$ws.Range("N3").Value = 2.343559
$ws.Range("O3").Value = 0.4884647274191578
$ws.Range("P3").Value = 0.4884647274191579
$ws.Range("Q3").Value = 27.91870223063178
$ws.Range("R3").Value = 251.268320075686
$ws.Range("S3").Value = 0.009523740442699885
$ws.Range("T3").Value = 0.009523740442699885

# Row 4
$ws.Range("H4").Value = 5067.86792
$ws.Range("I4").Value = 0.9215900675332435
$ws.Range("J4").Value = 0.9215900675332435
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8180823333333334
$ws.Range("N4").Value = 2.454247
$ws.Range("O4").Value = 0.5115352725808422
$ws.Range("P4").Value = 0.5115352725808422
$ws.Range("Q4").Value = 1381.977737672916
$ws.Range("R4").Value = 12437.79963905624
$ws.Range("S4").Value = 0.4714258264034145
$ws.Range("T4").Value = 0.4714258264034145

# Row 5
$ws.Range("H5").Value = 5067.86792
$ws.Range("I5").Value = 0.9215900675332435
$ws.Range("J5").Value = 0.9215900675332435
$ws.Range("M5").Value = 0.7811863333333333
$ws.Range("N5").Value = 2.343559
$ws.Range("O5").Value = 0.4884647274191578
$ws.Range("P5").Value = 0.4884647274191579
$ws.Range("Q5").Value = 1319.649719414142
$ws.Range("R5").Value = 11876.84747472728
$ws.Range("S5").Value = 0.450164241129829
$ws.Range("T5").Value = 0.4501642411298291

# Row 6
$ws.Range("G6").Value = 93.641553
$ws.Range("H6").Value = 280.924659
$ws.Range("I6").Value = 0.05108605424341119
$ws.Range("J6").Value = 0.05108605424341119
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8180823333333334
$ws.Range("N6").Value = 2.454247
$ws.Range("O6").Value = 0.5115352725808422
$ws.Range("P6").Value = 0.5115352725808422
$ws.Range("Q6").Value = 76.606500175197
$ws.Range("R6").Value = 689.4585015767731
$ws.Range("S6").Value = 0.02613231868248303
$ws.Range("T6").Value = 0.02613231868248303

# Row 7
$ws.Range("G7").Value = 93.641553
$ws.Range("H7").Value = 280.924659
$ws.Range("I7").Value = 0.05108605424341119
$ws.Range("J7").Value = 0.05108605424341119
$ws.Range("M7").Value = 0.7811863333333333
$ws.Range("N7").Value = 2.343559
$ws.Range("O7").Value = 0.4884647274191578
$ws.Range("P7").Value = 0.4884647274191579
$ws.Range("Q7").Value = 73.15150143570899
$ws.Range("R7").Value = 658.363512921381
$ws.Range("S7").Value = 0.02495373556092816
$ws.Range("T7").Value = 0.02495373556092816

# Row 8
$ws.Range("G8").Value = 14.34625366666667
$ws.Range("H8").Value = 43.038761
$ws.Range("I8").Value = 0.007826584134129748
$ws.Range("J8").Value = 0.007826584134129748
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8180823333333334
$ws.Range("N8").Value = 2.454247
$ws.Range("O8").Value = 0.5115352725808422
$ws.Range("P8").Value = 0.5115352725808422
$ws.Range("Q8").Value = 11.73641667421856
$ws.Range("R8").Value = 105.627750067967
$ws.Range("S8").Value = 0.004003573848428955
$ws.Range("T8").Value = 0.004003573848428955

# Row 9
$ws.Range("G9").Value = 14.34625366666667
$ws.Range("H9").Value = 43.038761
$ws.Range("I9").Value = 0.007826584134129748
$ws.Range("J9").Value = 0.007826584134129748
$ws.Range("M9").Value = 0.7811863333333333
$ws.Range("N9").Value = 2.343559
$ws.Range("O9").Value = 0.4884647274191578
$ws.Range("P9").Value = 0.4884647274191579
$ws.Range("Q9").Value = 11.20709729893322
$ws.Range("R9").Value = 100.863875690399
$ws.Range("S9").Value = 0.003823010285700793
$ws.Range("T9").Value = 0.003823010285700793
